# Update a couple runs as complete (current)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_runs")

# Row 34 (2050_TM152_DBP_NoProject_01): clear status (was "current")
$ws.Range("H34").Value = ""

# Row 35 (2050_TM152_DBP_NoProject_03): fill in urbansim_path/urbansim_runid
# and flip status from "running" to "current"
$ws.Range("F35").Value = """Blueprint Plus Crossing (s23)\v1.5.2"""
$ws.Range("G35").Value = "run72"
$ws.Range("H35").Value = "current"

# Row 40 (2050_TM152_DBP_PlusCrossing_02): clear status (was "current")
$ws.Range("H40").Value = ""

# Row 41 (2050_TM152_DBP_PlusCrossing_03): flip status from "running" to "current"
$ws.Range("H41").Value = "current"

# Update the active selection to match the saved view state
$ws.Range("G43").Select()
